{"js": "// Locate the field (the Word \"complex field\" that holds the M2Doc\n// expression `m:self.`) inside the document body, delete it, and\n// replace it with plain-text runs spelling out the same expression\n// wrapped in braces: \"{m:self.\" + \"\" + \"\" + \"}\" \u2014 mirroring the\n// TokenIteratorFieldRewriterSplit output (a field rewritten/split\n// into several plain runs instead of staying a real Word field).\n\nconst body = context.document.body;\n\nconst fields = body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length === 0) {\n  throw new Error(\"Expected at least one field in the document body.\");\n}\n\nconst field = fields.items[0];\n\n// Find the paragraph that contains the field, so we can re-insert the\n// replacement runs at exactly the spot the field used to occupy\n// (without hard-coding a paragraph index).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.fields.load(\"items\");\n}\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.fields.items.length > 0) {\n    targetParagraph = p;\n    break;\n  }\n}\nif (targetParagraph === null) {\n  throw new Error(\"Could not find the paragraph containing the field.\");\n}\n\n// Delete the field (removes begin/instrText/separate/result/end runs\n// as one unit) \u2014 this collapses to an empty cursor position where the\n// field used to start.\nfield.delete();\nawait context.sync();\n\n// In this document the field is the very first thing in its\n// paragraph, so the paragraph's \"Start\" is exactly where the field\n// used to begin.\nconst target = targetParagraph.getRange(\"Start\");\n\n// Four runs: \"{m:self.\" / \"\" / \"\" / \"}\" \u2014 each its own <w:r> with an\n// empty <w:rPr/>, matching the split performed by\n// TokenIteratorFieldRewriterSplit.\nconst runsXml =\n  \"<w:r><w:rPr/><w:t>{m:self.</w:t></w:r>\" +\n  \"<w:r><w:rPr/><w:t/></w:r>\" +\n  \"<w:r><w:rPr/><w:t/></w:r>\" +\n  \"<w:r><w:rPr/><w:t>}</w:t></w:r>\";\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  runsXml +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Replace the Word field holding the M2Doc expression \"m:self.\" with\n# plain-text runs spelling out the same expression wrapped in braces:\n# \"{m:self.\" + \"\" + \"\" + \"}\" (four runs, each with its own empty\n# <w:rPr/>), mirroring the split performed by\n# TokenIteratorFieldRewriterSplit (a real Word field rewritten into\n# several plain runs instead of staying a field).\n\n$d = $word.ActiveDocument\n\n$fields = $d.Fields\nif ($fields.Count -lt 1) {\n    throw \"Expected at least one field in the document.\"\n}\n$field = $fields.Item(1)\n\n# Find the character offset where the field starts, then locate the\n# paragraph that contains it (so we don't have to hard-code a\n# paragraph index).\n$fieldStart = $field.Code.Start\n\n$paragraphs = $d.Paragraphs\n$targetParaStart = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $p = $paragraphs.Item($i)\n    if (($fieldStart -ge $p.Range.Start) -and ($fieldStart -lt $p.Range.End)) {\n        $targetParaStart = $p.Range.Start\n        break\n    }\n}\nif ($targetParaStart -lt 0) {\n    throw \"Could not find the paragraph containing the field.\"\n}\n\n# Delete the whole field (begin/instrText/separate/result/end runs are\n# removed as a single unit), then insert the replacement runs exactly\n# where the field used to start.\n$field.Delete()\n\n$insertPoint = $d.Range($targetParaStart, $targetParaStart)\n\n$runsXml = \"<w:r><w:rPr/><w:t>{m:self.</w:t></w:r>\" +\n           \"<w:r><w:rPr/><w:t/></w:r>\" +\n           \"<w:r><w:rPr/><w:t/></w:r>\" +\n           \"<w:r><w:rPr/><w:t>}</w:t></w:r>\"\n\n$xml = \"<?xml version='1.0' encoding='UTF-8' standalone='yes'?>\" +\n       \"<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>\" +\n       \"<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>\" +\n       \"<pkg:xmlData>\" +\n       \"<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n       \"<w:body><w:p>$runsXml</w:p></w:body></w:document>\" +\n       \"</pkg:xmlData></pkg:part></pkg:package>\"\n\n$insertPoint.InsertXML($xml)\n"}
